# Adding the changes we made on may 9th
# Insert 7 new rows of accelerometer data right after the header row (row 1),
# pushing the existing 20 data rows down to rows 9-28, then append 3 more
# new rows at the end (rows 29-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 7 rows at the top of the data (before current row 2) ---
$ws.Range("A2:A8").EntireRow.Insert()
$ws.Range("A2:C8").ClearFormats()

$topRows = @(
    @(-1.765928411483764, 1.167147111892701, 0.571823143959045),
    @(-1.709035491943359, 1.218545722961426, 0.4037320613861088),
    @(-1.738995742797851, 1.102185392379761, 0.6806826651096344),
    @(-1.775501251220703, 1.23273515701294, 0.5835052132606506),
    @(-1.910298776626587, 1.120057487487793, 0.7890581786632539),
    @(-1.610696506500243, 1.25785665512085, 0.8155304193496704),
    @(-2.127863931655885, 1.3188805103302, 0.8300568103790283)
)

$r = 2
foreach ($row in $topRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Append 3 new rows at the bottom (rows 29-31) ---
$bottomRows = @(
    @(-0.8792205810546895, 1.9060809135437, 0.1410199522972106),
    @(-0.9211750030517577, 1.817085385322571, 0.09843596816062908),
    @(-0.8201595306396483, 1.820924615859986, 0.07434962689876565)
)

$r = 29
foreach ($row in $bottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
